# feat: add 2022-Q3 data
#
# Target layout after edit:
#   Sheets: 总计, 2022-Q3, 2022-Q2, 2022-Q1  (new "2022-Q3" sheet inserted
#   right after "总计" and before the existing "2022-Q2" sheet)
#
#   "总计" sheet gains a new row for 2022-Q3 (10 holdings, 0.76 value),
#   pushing the existing 2022-Q2 / 2022-Q1 rows down by one.
#
#   The new "2022-Q3" sheet carries the same shape/format as the existing
#   "2022-Q2" sheet (fund code / name / scale / position / ratio / value /
#   rank columns) but with 10 funds of its own data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Build the new "2022-Q3" worksheet by duplicating "2022-Q2" (so it
#    inherits identical column layout + cell formatting), inserting it
#    before "2022-Q2" in the tab order.
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2)
$q3 = $wb.Worksheets.Item("2022-Q2 (2)")
$q3.Name = "2022-Q3"

# Extend formatting from the last existing data row (row 7) down through
# row 11 so the four extra rows share the same look (style s=2 on col A)
# as the rest of the table before we pour in values.
$q3.Range("A7:H7").Copy()
$q3.Range("A8:H11").PasteSpecial(-4122)

# Fund-code / name / scale / position / ratio / value columns (B:G) are
# stored as text in this workbook (codes like "000594" must keep their
# leading zeros), so force the whole block to Text before writing values.
$q3.Range("B2:G11").NumberFormat = "@"

# Row 2
$q3.Range("A2").Value = 0
$q3.Range("B2").Value = "160926"
$q3.Range("C2").Value = "大成创业板两年定期开放混合A"
$q3.Range("D2").Value = "7.68"
$q3.Range("E2").Value = "64.14"
$q3.Range("F2").Value = "2.67"
$q3.Range("G2").Value = "0.2051"
$q3.Range("H2").Value = 8

# Row 3
$q3.Range("A3").Value = 1
$q3.Range("B3").Value = "000594"
$q3.Range("C3").Value = "大摩进取优选股票"
$q3.Range("D3").Value = "5.38"
$q3.Range("E3").Value = "91.94"
$q3.Range("F3").Value = "3.48"
$q3.Range("G3").Value = "0.1872"
$q3.Range("H3").Value = 7

# Row 4
$q3.Range("A4").Value = 2
$q3.Range("B4").Value = "001825"
$q3.Range("C4").Value = "建信中国制造2025股票A"
$q3.Range("D4").Value = "3.96"
$q3.Range("E4").Value = "86.05"
$q3.Range("F4").Value = "3.02"
$q3.Range("G4").Value = "0.1196"
$q3.Range("H4").Value = 8

# Row 5
$q3.Range("A5").Value = 3
$q3.Range("B5").Value = "009798"
$q3.Range("C5").Value = "大成创业板两年定期开放混合C"
$q3.Range("D5").Value = "2.71"
$q3.Range("E5").Value = "64.14"
$q3.Range("F5").Value = "2.67"
$q3.Range("G5").Value = "0.0724"
$q3.Range("H5").Value = 8

# Row 6
$q3.Range("A6").Value = 4
$q3.Range("B6").Value = "002707"
$q3.Range("C6").Value = "摩根士丹利华鑫科技领先灵活配置混合A"
$q3.Range("D6").Value = "1.75"
$q3.Range("E6").Value = "94.13"
$q3.Range("F6").Value = "3.22"
$q3.Range("G6").Value = "0.0564"
$q3.Range("H6").Value = 7

# Row 7
$q3.Range("A7").Value = 5
$q3.Range("B7").Value = "014380"
$q3.Range("C7").Value = "建信中国制造2025股票C"
$q3.Range("D7").Value = "1.69"
$q3.Range("E7").Value = "86.05"
$q3.Range("F7").Value = "3.02"
$q3.Range("G7").Value = "0.0510"
$q3.Range("H7").Value = 8

# Row 8
$q3.Range("A8").Value = 6
$q3.Range("B8").Value = "233011"
$q3.Range("C8").Value = "大摩主题优选混合"
$q3.Range("D8").Value = "1.34"
$q3.Range("E8").Value = "90.66"
$q3.Range("F8").Value = "3.46"
$q3.Range("G8").Value = "0.0464"
$q3.Range("H8").Value = 7

# Row 9
$q3.Range("A9").Value = 7
$q3.Range("B9").Value = "004223"
$q3.Range("C9").Value = "金信多策略精选灵活配置混合"
$q3.Range("D9").Value = "0.32"
$q3.Range("E9").Value = "92.79"
$q3.Range("F9").Value = "4.83"
$q3.Range("G9").Value = "0.0155"
$q3.Range("H9").Value = 8

# Row 10
$q3.Range("A10").Value = 8
$q3.Range("B10").Value = "350007"
$q3.Range("C10").Value = "天治趋势精选混合"
$q3.Range("D10").Value = "0.37"
$q3.Range("E10").Value = "82.33"
$q3.Range("F10").Value = "1.97"
$q3.Range("G10").Value = "0.0073"
$q3.Range("H10").Value = 10

# Row 11
$q3.Range("A11").Value = 9
$q3.Range("B11").Value = "014871"
$q3.Range("C11").Value = "摩根士丹利华鑫科技领先灵活配置混合C"
$q3.Range("D11").Value = "0.08"
$q3.Range("E11").Value = "94.13"
$q3.Range("F11").Value = "3.22"
$q3.Range("G11").Value = "0.0026"
$q3.Range("H11").Value = 7

# ---------------------------------------------------------------------
# 2. Update the "总计" (totals) sheet: insert a 2022-Q3 row at the top of
#    the data (row 2), pushing 2022-Q2 / 2022-Q1 down by one row.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Give the brand-new row 4 the same look (style) as the existing index
# cells in column A before writing into it.
$total.Range("A3").Copy()
$total.Range("A4").PasteSpecial(-4122)

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2022-Q1"
$total.Range("C4").Value = 3
$total.Range("D4").Value = 0.51

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q2"
$total.Range("C3").Value = 6
$total.Range("D3").Value = 0.83

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 10
$total.Range("D2").Value = 0.76

# ---------------------------------------------------------------------
# 3. Restore "2022-Q1" as the active / selected sheet (it was the tab
#    selected in the original workbook).
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("2022-Q1")
$q1.Activate()
